$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: shift the existing quarter rows (2..8) down by
#    one (3..9) and write the new 2022-Q4 row into row 2.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

for ($r = 8; $r -ge 2; $r--) {
    $src = $summary.Range("A" + $r + ":D" + $r)
    $dst = $summary.Range("A" + ($r + 1) + ":D" + ($r + 1))
    $src.Copy($dst)
}

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 0.48

# ---------------------------------------------------------------------------
# 2. Insert a brand-new worksheet named "2022-Q4" right before "2022-Q3" and
#    populate it with the quarterly fund-holding detail.
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"

# Match the page-setup margins used by every other quarter sheet.
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# Re-fetch "2022-Q3" now that the sheet collection has shifted - it is the
# style donor for the header / index-column formatting.
$q3 = $wb.Worksheets.Item("2022-Q3")

# Header row (bold / centered / bordered style copied from the sibling sheet)
$q3.Range("B1:H1").Copy($q4.Range("B1:H1"))
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Index column (A2:A8) uses the same styled-index look as every other sheet.
for ($r = 2; $r -le 8; $r++) {
    $q3.Range("A2").Copy($q4.Range("A" + $r))
    $q4.Range("A" + $r).Value = ($r - 2)
}

# Force the text columns (B..G) to Text so that fund codes keep their
# leading zeros and the percentages keep their trailing zeros, matching the
# inline-string data in the other quarter sheets.
$q4.Range("B2:G8").NumberFormat = "@"

$rows = @(
    @("014158", "博时浦惠一年持有期混合A", "3.84", "49.14", "3.96", "0.1521", 2),
    @("002095", "博时新收益灵活配置混合A", "4.96", "88.69", "3.03", "0.1503", 10),
    @("002096", "博时新收益灵活配置混合C", "4.21", "88.69", "3.03", "0.1276", 10),
    @("004189", "华商消费行业股票", "0.73", "80.96", "2.87", "0.0210", 9),
    @("014159", "博时浦惠一年持有期混合C", "0.37", "49.14", "3.96", "0.0147", 2),
    @("010663", "长江均衡成长混合A", "0.21", "86.40", "5.47", "0.0115", 1),
    @("010664", "长江均衡成长混合C", "0.05", "86.40", "5.47", "0.0027", 1)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $row = $rows[$i]
    $q4.Range("B" + $r).Value = $row[0]
    $q4.Range("C" + $r).Value = $row[1]
    $q4.Range("D" + $r).Value = $row[2]
    $q4.Range("E" + $r).Value = $row[3]
    $q4.Range("F" + $r).Value = $row[4]
    $q4.Range("G" + $r).Value = $row[5]
    $q4.Range("H" + $r).Value = $row[6]
}

# Restore the workbook's original active sheet ("总计") so the edit doesn't
# change the view state beyond what's needed for the new data.
$summary.Activate()

